$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlLeft = -4131
$xlCenter = -4108

# Fill in the missing Issue Type / Resolved? for the existing row 8
$ws.Range("B8").Value = "Bugfix"
$ws.Range("B8").HorizontalAlignment = $xlCenter
$ws.Range("C8").Value = "Yes"
$ws.Range("C8").HorizontalAlignment = $xlCenter

# New row 9: feature request, no detail/comment
$ws.Range("A9").Value = "Need to support 632 by 1030"
$ws.Range("A9").HorizontalAlignment = $xlLeft
$ws.Range("B9").Value = "Feature Request"
$ws.Range("B9").HorizontalAlignment = $xlCenter

# New row 10: bugfix, no detail/comment
$ws.Range("A10").Value = "After Deploy 1 unit, bot just sits there"
$ws.Range("A10").HorizontalAlignment = $xlLeft
$ws.Range("B10").Value = "Bugfix"
$ws.Range("B10").HorizontalAlignment = $xlCenter

# New row 11: feature request with detail
$ws.Range("A11").Value = "Bot needs to retry failed quests"
$ws.Range("A11").HorizontalAlignment = $xlLeft
$ws.Range("B11").Value = "Feature Request"
$ws.Range("B11").HorizontalAlignment = $xlCenter
$ws.Range("D11").Value = "Edit the BotAllQuests.ahk"
$ws.Range("D11").HorizontalAlignment = $xlLeft

# Update selection to match the recorded state after editing
$ws.Range("D10").Select()
